## Changed to manage multiple sheets.
##
## 1) "My sheet" gets some light formatting:
##      - A2 (rahul)  -> bold
##      - A3 (priya)  -> blue, 20pt
##      - A5           -> new empty cell boxed with a double border
## 2) A second worksheet "My sheet 2" is added right after "My sheet",
##    holding a single note in C1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- formatting on the existing sheet -------------------------------------

# "rahul" -> bold
$ws1.Range("A2").Font.Bold = $true

# "priya" -> blue, size 20
$priyaFont = $ws1.Range("A3").Font
$priyaFont.Color = 16711680   # blue (OLE BGR encoding of RGB 0,0,255)
$priyaFont.Size = 20

# A5 -> empty cell surrounded by a double border
$ws1.Range("A5").Borders.LineStyle = -4119   # xlDouble

# --- add the second worksheet ---------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "My sheet 2"
$ws2.Range("C1").Value = "writing ;)"

# Leave the first sheet active / A1 selected, same as the original workbook.
$ws1.Activate()
[void]$ws1.Range("A1").Select()

Write-Output "done"
